$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The LED resistor (row 28) value changed from 130Ω to 120Ω, with its
# package updated from the 0603 footprint to the 0402 footprint to match.
# A leading apostrophe is used so the text-literal (quote-prefixed) cell
# style already applied to these cells is preserved instead of being reset.
$ws.Range("B28").Value = "'RESISTOR, 120Ω"
$ws.Range("D28").Value = "'Ultiboard\R0402"

$ws.Range("B28").Select()
